# 13.1.1 "Number of deaths attributed to disasters" sheet: add the 2022
# year column (S) after the existing 2021 column (R), with a data point
# for every row that already has per-year figures (rows 3-34), and move
# the active selection to T4 (the cell immediately to the right of the
# newly added column's header), matching the author's saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at S. Excel extends the formatting of the
# neighbouring column (R) onto the freshly inserted one automatically,
# which reproduces the per-row cell styles (borders/fonts/number formats)
# without fabricating any new style entries.
[void]$ws.Columns("S").Insert()

# Row 4 holds the year headers.
$ws.Range("S4").Value = 2022

# Row 5: "Кыргыз Республикасы" / "Кыргызская Республика" / "Kyrgyz Republic" totals.
$ws.Range("S5").Value = 135
$ws.Range("S6").Value = 99
$ws.Range("S7").Value = 36

# Row 8: "Баткен облусу" / "Баткенская область" / "Batken oblast".
$ws.Range("S8").Value = 97
$ws.Range("S9").Value = 80
$ws.Range("S10").Value = 17

# Row 11: "Жалал-Абад облусу" / "Джалал-Абадская область" / "Djalal-Abad oblast".
$ws.Range("S11").Value = 17
$ws.Range("S12").Value = 11
$ws.Range("S13").Value = 6

# Row 14: "Ысык-Көл облусу" / "Иссык-Кульская область" / "Ysyk-Kul oblast".
$ws.Range("S14").Value = 5
$ws.Range("S15").Value = 3
$ws.Range("S16").Value = 2

# Row 17: "Нарын облусу" / "Нарынская область" / "Naryn oblast" - no data ("-").
$ws.Range("S17").Value = "-"
$ws.Range("S18").Value = "-"
$ws.Range("S19").Value = "-"

# Row 20: "Ош облусу" / "Ошская область" / "Osh oblast".
$ws.Range("S20").Value = 6
$ws.Range("S21").Value = 1
$ws.Range("S22").Value = 5

# Row 23: "Талас облусу" / "Таласская область" / "Talas oblast" - no data ("-").
$ws.Range("S23").Value = "-"
$ws.Range("S24").Value = "-"
$ws.Range("S25").Value = "-"

# Row 26: "Чүй облусу" / "Чуйская область" / "Chui oblast".
$ws.Range("S26").Value = 10
$ws.Range("S27").Value = 4
$ws.Range("S28").Value = 6

# Row 29: "Бишкек ш." / "г.Бишкек" / "Bishkek city" - no data ("-").
$ws.Range("S29").Value = "-"
$ws.Range("S30").Value = "-"
$ws.Range("S31").Value = "-"

# Row 32: "Ош ш." / "г.Ош" / "Osh city" - no data ("-") for all three sub-rows.
$ws.Range("S32").Value = "-"
$ws.Range("S33").Value = "-"
$ws.Range("S34").Value = "-"

# Match the author's final selection/view state.
[void]$ws.Range("T4").Select()
